$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.243.36"
$ws.Range("D3").Value = "2.227.04"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'244.68"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'73.68"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").Value = "'42.58"
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("D11").Value = "'0.0966"
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("D12").Value = "'7.13"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'0.851"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "2.234.23"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "42.119.33"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "'0.0000112"
$ws.Range("E18").Value = "  +16.43%  "
$ws.Range("D19").Value = "'6.16"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "'72.11"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "'9.94"
$ws.Range("E21").Value = "  +38.32%  "
$ws.Range("D22").Value = "'231.26"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "'11.88"
$ws.Range("E24").Value = "  +8.05%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").Value = "'166.92"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'21.07"
$ws.Range("E30").Value = "  +3.68%  "
$ws.Range("D31").Value = "'5.75"
$ws.Range("E31").Value = "  +19.59%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "'0.118"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("D35").Value = "'29.46"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").Value = "'13.02"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "'5.62"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").Value = "'62.96"
$ws.Range("E41").Value = "  +5.86%  "
$ws.Range("D42").Value = "'0.201"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "'8.82"
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("D44").Value = "'105.29"
$ws.Range("E44").Value = "  -5.01%  "
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("E47").Value = "  +7.26%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "'1.17"
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").Value = "'4.06"
$ws.Range("E51").Value = "  +0.77%  "

# Reset style on cells that needed an apostrophe-prefix to stay text
# (so numeric-looking strings like "244.68" are not auto-converted to numbers),
# to avoid leaving a stray quote-prefix style behind.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
